$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude = $wb.Worksheets.Item("Include #0")

# URL value (ValueSet URL changed from NPHCDA to HL7 R4 administrative-gender)
$wsMeta.Range("B2").Value = "http://hl7.org/fhir/R4/valueset-administrative-gender"

# Status value
$wsMeta.Range("B6").Value = "draft"

# Date value
$wsMeta.Range("B8").Value = "2025-06-25T06:29:04+01:00"

# Description value
$wsMeta.Range("B13").Value = "Subset of HL7 administrative-gender limited to 'male' and 'female'."

# System URI value on the Include sheet now matches the (new) ValueSet URL
$wsInclude.Range("B4").Value = "http://hl7.org/fhir/R4/valueset-administrative-gender"
